# "Actualizar" refresh macro: shifts the rolling availability-check
# timestamps down one slot and stamps the newest slot (rows 2-15) with the
# current run time (02-16-2021 12:07:57).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44243.50540237159
$ws.Range("D16:D29").Value = 44243.48416751157
$ws.Range("D30:D43").Value = 44243.46292204861
